# Rename the header labels on row 1 from the uppercase abbreviations
# (RIQ, ABUN, H') to the new lowercase labels (riq, abund, shan).
# Because the old labels become unused elsewhere in the sheet, the
# workbook's shared-strings table drops them and appends the new ones,
# which naturally re-numbers the remaining shared strings (HORTA, MATA,
# PASTO) to the front of the table - matching the target layout without
# needing to touch the data rows themselves.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "riq"
$ws.Range("C1").Value = "abund"
$ws.Range("D1").Value = "shan"

# Update the active selection to D2 (was the whole-sheet selection before).
$ws.Range("D2").Select()
